# "Correção metricas de CK & LK"
#
# 1. Decrement the WMC (col B) and CS (col F) metric values by 1 for every
#    data row (3-21) on the "Metricas CK & LK" sheet.
# 2. Make "Metricas CK & LK" the active sheet/tab (was "CC McCabe"), with
#    the selection moved from I21 to F21.
# 3. On "CC McCabe", the selection stays C10 but the view scrolls so that
#    A47 is the top-left visible cell (and it is no longer the active tab).

$wb = $excel.ActiveWorkbook

$wsCk = $wb.Worksheets.Item("Metricas CK & LK")
$wsCc = $wb.Worksheets.Item("CC McCabe")

# --- Correct the CK/LK metric values (WMC = col B, CS = col F) ---------
for ($row = 3; $row -le 21; $row++) {
    $wsCk.Cells.Item($row, 2).Value = $wsCk.Cells.Item($row, 2).Value2 - 1
    $wsCk.Cells.Item($row, 6).Value = $wsCk.Cells.Item($row, 6).Value2 - 1
}

# --- Fix up view/selection state ---------------------------------------
# CC McCabe keeps its C10 selection, but scroll the view so A47 sits at
# the top-left before we leave the sheet (and before it loses tab focus).
$wsCc.Activate()
$wsCc.Range("A47").Select()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$wsCc.Range("C10").Select()

# Activate "Metricas CK & LK" last so it becomes the active tab/sheet, and
# select F21 on it.
$wsCk.Activate()
$wsCk.Range("F21").Select()
